# Generate Report for Handoff
# -------------------------------------------------------------------
# The b.md file has now been handed off for zh-cn and de-de.  This
# refreshes the "Overview" summary sheet plus the two per-locale
# status sheets (zh-cn, de-de) for the b.md row with:
#   - Status => "Ready for handoff"
#   - a new "Latest Handoff File" (the freshly generated xlf)
#   - a new "Latest Handoff Datetime"
#   - an "Error Detail" noting the handback file version is stale
# It also widens the "Error Detail" column on both locale sheets so
# the longer message is readable.
# -------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$status = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/e0f0a154746cf7bb58b4775775c7405801b460ca/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/175d2352a915787f1b9b0c1f3c639821dc22df30/e2e/b.md."

# ---------------- Overview sheet : row 3 is b.md ----------------
$overview.Range("E3").Value = $status
$overview.Range("F3").Value = $status
$overview.Range("G3").Value = "2016-08-02 09:43:15"

# ---------------- zh-cn sheet : row 3 is b.md --------------------
$zhcn.Range("C3").Value = $status
$zhcn.Range("F3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("G3").Value = "2016-08-02 09:43:05"
$zhcn.Range("O3").Value = $errorDetail
$zhcn.Columns.Item(15).ColumnWidth = 39.17

# ---------------- de-de sheet : row 3 is b.md --------------------
$dede.Range("C3").Value = $status
$dede.Range("F3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("G3").Value = "2016-08-02 09:43:15"
$dede.Range("O3").Value = $errorDetail
$dede.Columns.Item(15).ColumnWidth = 39.17
